$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# 1) Rows 235-239 (match ids / odds for 2024-04-15 Liga I round):
#    The underlying source data got re-sorted. Column A (sequence
#    number) stays fixed per row; every other column (B..AC) moves
#    between rows according to the mapping below.
#    new row -> old row that supplies its B..AC content
# -----------------------------------------------------------------
$rowMap = @{
    235 = 239
    236 = 238
    237 = 235
    238 = 237
    239 = 236
}

$firstCol = 2   # column B
$lastCol  = 29  # column AC

# Snapshot the current (pre-edit) contents of every source row first,
# so overwriting one row doesn't clobber data another row still needs.
$snapshot = @{}
foreach ($srcRow in $rowMap.Values) {
    if (-not $snapshot.ContainsKey($srcRow)) {
        $rowVals = @{}
        for ($col = $firstCol; $col -le $lastCol; $col++) {
            $rowVals[$col] = $ws.Cells.Item($srcRow, $col).Value2
        }
        $snapshot[$srcRow] = $rowVals
    }
}

foreach ($destRow in $rowMap.Keys) {
    $srcRow = $rowMap[$destRow]
    $rowVals = $snapshot[$srcRow]
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $ws.Cells.Item($destRow, $col).Value = $rowVals[$col]
    }
}

# -----------------------------------------------------------------
# 2) Rows 267, 268, 270, 272, 273, 274: refreshed closing odds
#    (oddH/oddD/oddA/Ah/oddAHH/oddAHA/oddAHOver/oddAHUnder columns).
# -----------------------------------------------------------------
$oddsUpdates = @{
    267 = @{ N = 2.375; P = 3;    Q = -0.25; R = 2.125; S = 1.75;  U = 1.925; V = 1.925 }
    268 = @{ N = 1.85;  P = 4;               R = 1.875; S = 1.975; U = 2;     V = 1.85  }
    270 = @{                                                       U = 1.85;  V = 2     }
    272 = @{            O = 3.2;                                   U = 1.925; V = 1.925 }
    273 = @{ N = 2.5;   P = 2.8; Q = 0;      R = 1.8;   S = 2.05 }
    274 = @{ N = 1.85;  O = 3.3; P = 4;                S = 1.925 }
}

foreach ($r in $oddsUpdates.Keys) {
    $cols = $oddsUpdates[$r]
    foreach ($colLetter in $cols.Keys) {
        $addr = "$colLetter$r"
        $ws.Range($addr).Value = $cols[$colLetter]
    }
}
